# Add a new "Italy" Test Data worksheet, mirroring the layout of the existing
# "Germany" sheet, and update the sheet selections/active tab to match.

$wb = $excel.ActiveWorkbook

$wsGermany  = $wb.Worksheets.Item("Germany")
$wsSlovakia = $wb.Worksheets.Item("Slovakia")

# Create the new sheet by copying "Germany" (same column widths / cell styles /
# merged cells) and place it after the last sheet ("Slovakia").
$wsGermany.Copy($null, $wsSlovakia) | Out-Null
$wsItaly = $wb.Worksheets.Item($wsSlovakia.Index + 1)
$wsItaly.Name = "Italy"

# Fill in the Italy-specific values (market name + Jira ticket references).
$wsItaly.Range("B2").Value = "Italy Market"
$wsItaly.Range("B4").Value = "NGC-3145/T2454/T2453/T2452/T2455/NGC-3145/T2446"

# Update view/selection state:
#  - Germany now has its whole grid selected.
#  - Slovakia's selection moves to B9, and it is no longer the active tab.
#  - The new Italy sheet becomes the active tab, with row 12 selected.
$wsGermany.Cells.Select() | Out-Null
$wsSlovakia.Range("B9").Select() | Out-Null
$wsItaly.Rows.Item(12).Select() | Out-Null
